$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualización lista de materiales: update quantities in column B
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("B6").Value = 6
$ws.Range("B12").Value = 6
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 12

$ws.Range("B13").Select()
